$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1; existing rows 1-20 shift down to 2-21.
$ws.Rows("1:1").Insert()

# Populate the new header row.
$ws.Range("A1").Value = "Reference"
$ws.Range("B1").Value = "UNAM"

# Update the sheet view selection to B1 (was E1).
$ws.Range("B1").Select()

# The hyperlinks collection does not auto-shift with the row insert, so
# remove the old ones and re-create them at their new (shifted) locations.
$ws.Range("B18").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B19"), "https://en.wikipedia.org/wiki/Michoac%C3%A1n", "", "Michoacán", "https://en.wikipedia.org/wiki/Michoac%C3%A1n")
$ws.Range("B19").Value = "Michoacán"

$ws.Hyperlinks.Add($ws.Range("B21"), "https://en.wikipedia.org/wiki/Veracruz", "", "Veracruz", "https://en.wikipedia.org/wiki/Veracruz")
$ws.Range("B21").Value = "Veracruz"
